$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "46.198.70"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "2.605.66"
$ws.Range("E3").Value = "  +10.10%  "

$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "310.65"
$ws.Range("E5").Value = "  +3.95%  "

Set-TextValue $ws.Range("D6") "100.05"
$ws.Range("E6").Value = "  +3.18%  "

Set-TextValue $ws.Range("D7") "0.597"

$ws.Range("E8").Value = "  +0.00%  "

Set-TextValue $ws.Range("D9") "0.585"
$ws.Range("E9").Value = "  +16.72%  "

Set-TextValue $ws.Range("D10") "38.63"
$ws.Range("E10").Value = "  +14.24%  "

Set-TextValue $ws.Range("D11") "0.0840"
$ws.Range("E11").Value = "  +7.40%  "

Set-TextValue $ws.Range("D12") "8.35"
$ws.Range("E12").Value = "  +19.14%  "

$ws.Range("D13").Value = "3.000.37"
$ws.Range("E13").Value = "  +9.82%  "

$ws.Range("E14").Value = "  +1.63%  "

$ws.Range("D15").Value = "2.599.05"
$ws.Range("E15").Value = "  +9.59%  "

$ws.Range("E16").Value = "  +11.39%  "

Set-TextValue $ws.Range("D17") "14.88"
$ws.Range("E17").Value = "  +9.26%  "

$ws.Range("D18").Value = "46.394.69"
$ws.Range("E18").Value = "  +1.73%  "

$ws.Range("E19").Value = "  +7.86%  "

Set-TextValue $ws.Range("D20") "13.05"
$ws.Range("E20").Value = "  +3.40%  "

Set-TextValue $ws.Range("D21") "6.70"
$ws.Range("E21").Value = "  +11.79%  "

Set-TextValue $ws.Range("D22") "71.14"
$ws.Range("E22").Value = "  +6.76%  "

Set-TextValue $ws.Range("D23") "255.05"
$ws.Range("E23").Value = "  +5.44%  "

Set-TextValue $ws.Range("D24") "3.07"
$ws.Range("E24").Value = "  +11.52%  "

Set-TextValue $ws.Range("D25") "2.21"
$ws.Range("E25").Value = "  +17.22%  "

Set-TextValue $ws.Range("D26") "28.31"
$ws.Range("E26").Value = "  +35.81%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("E28").Value = "  +9.70%  "

Set-TextValue $ws.Range("D29") "39.72"
$ws.Range("E29").Value = "  +4.16%  "

Set-TextValue $ws.Range("D30") "2.27"
$ws.Range("E30").Value = "  +3.51%  "

Set-TextValue $ws.Range("D31") "6.12"
$ws.Range("E31").Value = "  +12.18%  "

$ws.Range("E32").Value = "  -0.87%  "

Set-TextValue $ws.Range("D33") "2.29"
$ws.Range("E33").Value = "  +21.33%  "

Set-TextValue $ws.Range("D34") "2.89"
$ws.Range("E34").Value = "  +5.13%  "

Set-TextValue $ws.Range("D35") "152.86"
$ws.Range("E35").Value = "  +4.44%  "

Set-TextValue $ws.Range("D36") "0.0831"
$ws.Range("E36").Value = "  +9.06%  "

$ws.Range("E37").Value = "  +4.80%  "

$ws.Range("E38").Value = "  +5.92%  "

Set-TextValue $ws.Range("D39") "16.50"
$ws.Range("E39").Value = "  +9.85%  "

Set-TextValue $ws.Range("D40") "4.19"
$ws.Range("E40").Value = "  +10.02%  "

Set-TextValue $ws.Range("D41") "3.60"
$ws.Range("E41").Value = "  +12.70%  "

Set-TextValue $ws.Range("D42") "21.69"
$ws.Range("E42").Value = "  +54.21%  "

Set-TextValue $ws.Range("D43") "0.0325"
$ws.Range("E43").Value = "  +10.57%  "

$ws.Range("D44").Value = "2.040.62"
$ws.Range("E44").Value = "  +5.16%  "

$ws.Range("E45").Value = "  -0.03%  "

Set-TextValue $ws.Range("D46") "91.09"
$ws.Range("E46").Value = "  -1.85%  "

Set-TextValue $ws.Range("D47") "9.26"
$ws.Range("E47").Value = "  +9.42%  "

Set-TextValue $ws.Range("D48") "109.57"
$ws.Range("E48").Value = "  +11.26%  "

$ws.Range("E49").Value = "  +1.12%  "

Set-TextValue $ws.Range("D50") "0.200"
$ws.Range("E50").Value = "  +10.07%  "

$ws.Range("D51").Value = "2.862.05"
$ws.Range("E51").Value = "  +9.98%  "
